$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had 4 columns of button/element locators (A:D) and
# two data rows. The edit trims the sheet down to a single column (the
# "input_KeyName" locator) with its row of (now blank) sample data.

# Drop columns B:D entirely - shifts the grid left so only column A remains,
# which also shrinks the worksheet's used range/dimension down to column A.
$ws.Range("B:D").Delete()

# Column A keeps its original locator width, but it now carries the value
# (and width) that used to live in column D.
$ws.Columns.Item(1).ColumnWidth = 14.1666666666667

# Row 1: the single remaining header cell becomes the "input_KeyName" locator.
$ws.Range("A1").Value = "input_KeyName"

# Row 2: the matching sample value is now blank - clear it back to an empty
# cell (re-applying the default style keeps an empty <c> placeholder for A2
# instead of dropping the cell/row entirely).
$ws.Range("A2").Value = ""
$ws.Range("A2").Style = "Normal"
